# Commit: "Tue, Mar 31, 2020  5:04:59 AM"
#
# The presentation's main theme (ppt/theme/theme1.xml, the "Integral"
# theme used by the slide master / all slides) is swapped for the plain
# default "Office Theme" that already shipped in this file as
# ppt/theme/theme2.xml (only used by the notes master).
#
# The two themes' font scheme and format scheme (fills/lines/effects) are
# byte-for-byte identical - the only real difference between them is the
# <a:clrScheme> color palette. So reproducing the edit is a matter of
# repointing the slide master's theme colors at the default Office
# palette.
#
# PowerPoint's ColorScheme object exposes the theme's twelve colors as a
# flat, 1-based list, in this fixed order:
#   1 dk1  2 lt1  3 dk2  4 lt2  5 accent1  6 accent2  7 accent3
#   8 accent4  9 accent5  10 accent6  11 hlink  12 folHlink
# dk1/lt1 (pure black/white) are already identical between the two
# themes, so only items 3-12 need to be updated.

function HexToRgbInt([string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return $r + ($g * 256) + ($b * 65536)
}

$p = $ppt.ActivePresentation
$colors = $p.SlideMaster.ColorScheme

$colors.Item(3).RGB  = HexToRgbInt "44546A"  # dk2
$colors.Item(4).RGB  = HexToRgbInt "E7E6E6"  # lt2
$colors.Item(5).RGB  = HexToRgbInt "5B9BD5"  # accent1
$colors.Item(6).RGB  = HexToRgbInt "ED7D31"  # accent2
$colors.Item(7).RGB  = HexToRgbInt "A5A5A5"  # accent3
$colors.Item(8).RGB  = HexToRgbInt "FFC000"  # accent4
$colors.Item(9).RGB  = HexToRgbInt "4472C4"  # accent5
$colors.Item(10).RGB = HexToRgbInt "70AD47"  # accent6
$colors.Item(11).RGB = HexToRgbInt "0563C1"  # hlink
$colors.Item(12).RGB = HexToRgbInt "954F72"  # folHlink
